# Append a new "2021年" data row (row 13) to Sheet1, mirroring the layout
# of the preceding rows (e.g. row 12): a styled year label in column A,
# numeric totals in most columns, and a handful of columns (D, J, M) that
# carry an explicit empty-text value rather than being truly blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 12's A-cell formatting (bold font, border, center/top alignment)
# onto A13 so the new year label reuses the existing header style.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 6

# D13/J13/M13 need to be empty TEXT cells (like D12/J12/M12) rather than
# blank numeric cells, so we write a lone "'" (Excel's text-prefix marker,
# which collapses to an empty string) and then re-apply the plain
# formatting from the corresponding row-12 cell.
$ws.Range("D13").Value = "'"
$ws.Range("D13").Style = $ws.Range("D12").Style

$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 6

$ws.Range("J13").Value = "'"
$ws.Range("J13").Style = $ws.Range("J12").Style

$ws.Range("K13").Value = 6
$ws.Range("L13").Value = 13

$ws.Range("M13").Value = "'"
$ws.Range("M13").Style = $ws.Range("M12").Style

$ws.Range("N13").Value = 19
$ws.Range("O13").Value = 63
